$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Rows where the "In-charge Full Name" (column B) changes from "Admin" to "Staff"
$ws.Range("B28").Value = "Staff"
$ws.Range("B29").Value = "Staff"
$ws.Range("B31").Value = "Staff"
$ws.Range("B32").Value = "Staff"

# Rows that get "Iteration 1" set as the Planned Code Iteration (column D)
$ws.Range("D9").Value = "Iteration 1"
$ws.Range("D12").Value = "Iteration 1"
$ws.Range("D18").Value = "Iteration 1"
$ws.Range("D19").Value = "Iteration 1"
$ws.Range("D20").Value = "Iteration 1"
$ws.Range("D27").Value = "Iteration 1"
$ws.Range("D28").Value = "Iteration 1"
$ws.Range("D29").Value = "Iteration 1"
$ws.Range("D30").Value = "Iteration 1"
$ws.Range("D31").Value = "Iteration 1"
$ws.Range("D32").Value = "Iteration 1"
$ws.Range("D40").Value = "Iteration 1"
